$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume(1h) updates for unchanged coin rows ---
# NumberFormat "@" (Text) is applied before assigning D/E values so that
# Excel does not silently reinterpret formatted numeric-looking strings
# (e.g. "58.97", "0.0215") as floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.848.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.087.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.97"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.391.68"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.24"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.084.27"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.759.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.40"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.00%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.70"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0632"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.47"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.30%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.70"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.466.56"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.50%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.55%  "

# --- Row 43/44: FTXToken and VeChain swap ranking positions ---
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.83%  "

# --- Row 51: RocketPoolETH replaced by MultiversX ---
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.38"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.22%  "
